$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column cells whose refreshed price text happens to look like a plain
# decimal number (e.g. "379.48"). Excel would otherwise auto-convert a
# bare numeric-looking Value to a number type, but the source keeps these
# as text (as with the thousand-dot values like "51.638.22"). Force the
# "Text" number format first so the assigned value is stored as a string.
$textForceRows = @(5,6,10,12,14,17,18,20,21,23,24,26,28,32,33,34,35,36,39,40,42,46,51)
foreach ($r in $textForceRows) {
    $ws.Range("D" + $r).NumberFormat = "@"
}

# Apply the refreshed price (D) and 1h volume change (E) figures.
$ws.Range("D2").Value = "51.638.22"
$ws.Range("D3").Value = "3.021.68"
$ws.Range("E3").Value = "  +2.18%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "379.48"
$ws.Range("E5").Value = "  -0.13%  "
$ws.Range("D6").Value = "102.45"
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("E7").Value = "  +0.42%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  +0.64%  "
$ws.Range("D10").Value = "36.83"
$ws.Range("E10").Value = "  +1.01%  "
$ws.Range("E11").Value = "  -0.20%  "
$ws.Range("D12").Value = "0.0863"
$ws.Range("E12").Value = "  +1.34%  "
$ws.Range("D13").Value = "3.502.88"
$ws.Range("E13").Value = "  +2.08%  "
$ws.Range("D14").Value = "18.48"
$ws.Range("E14").Value = "  +0.39%  "
$ws.Range("E15").Value = "  -0.38%  "
$ws.Range("D16").Value = "3.007.15"
$ws.Range("E16").Value = "  +1.81%  "
$ws.Range("D17").Value = "0.975"
$ws.Range("E17").Value = "  -3.59%  "
$ws.Range("D18").Value = "10.62"
$ws.Range("E18").Value = "  -14.30%  "
$ws.Range("D19").Value = "51.630.76"
$ws.Range("E19").Value = "  +1.02%  "
$ws.Range("D20").Value = "3.10"
$ws.Range("E20").Value = "  +1.03%  "
$ws.Range("D21").Value = "12.44"
$ws.Range("E21").Value = "  +0.29%  "
$ws.Range("E22").Value = "  +0.72%  "
$ws.Range("D23").Value = "70.03"
$ws.Range("E23").Value = "  +0.49%  "
$ws.Range("D24").Value = "267.10"
$ws.Range("E24").Value = "  -0.84%  "
$ws.Range("E25").Value = "  -6.20%  "
$ws.Range("D26").Value = "8.31"
$ws.Range("E26").Value = "  +4.10%  "
$ws.Range("E27").Value = "  +8.02%  "
$ws.Range("D28").Value = "0.172"
$ws.Range("E28").Value = "  +3.65%  "
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("E30").Value = "  +1.34%  "
$ws.Range("E31").Value = "  +0.54%  "
$ws.Range("D32").Value = "10.26"
$ws.Range("E32").Value = "  -2.64%  "
$ws.Range("D33").Value = "2.12"
$ws.Range("E33").Value = "  +2.68%  "
$ws.Range("D34").Value = "50.56"
$ws.Range("E34").Value = "  -1.13%  "
$ws.Range("D35").Value = "33.87"
$ws.Range("E35").Value = "  -1.14%  "
$ws.Range("D36").Value = "0.0452"
$ws.Range("E36").Value = "  +3.91%  "
$ws.Range("E37").Value = "  -0.11%  "
$ws.Range("E38").Value = "  +1.77%  "
$ws.Range("D39").Value = "0.293"
$ws.Range("E39").Value = "  +13.76%  "
$ws.Range("D40").Value = "16.97"
$ws.Range("E40").Value = "  +1.17%  "
$ws.Range("E41").Value = "  +1.27%  "
$ws.Range("D42").Value = "127.78"
$ws.Range("E42").Value = "  +5.72%  "
$ws.Range("E43").Value = "  -0.63%  "
$ws.Range("E44").Value = "  +1.20%  "
$ws.Range("E45").Value = "  +4.85%  "
$ws.Range("D46").Value = "21.58"
$ws.Range("E46").Value = "  -0.63%  "
$ws.Range("E47").Value = "  +2.38%  "
$ws.Range("E48").Value = "  +2.88%  "
$ws.Range("D49").Value = "2.025.73"
$ws.Range("E49").Value = "  -1.58%  "
$ws.Range("D50").Value = "3.319.39"
$ws.Range("E50").Value = "  +2.01%  "
$ws.Range("D51").Value = "0.0318"
$ws.Range("E51").Value = "  -1.57%  "

# Restore the default (Normal) cell style now that the text values are
# committed, so no stray "Text" number-format style lingers on the cells.
foreach ($r in $textForceRows) {
    $ws.Range("D" + $r).Style = "Normal"
}
